# MarsFramework/ExcelData/TestData.xlsx update
# - add Education, Language sheets (education tab fixes / dropdown lists)
# - add Certificate sheet (reads existing ISTQB/ANZTB values)
# - add blank placeholder sheets Sheet4-Sheet7
# - tweak remembered selections on SignIn / Profile
# - make Education the active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SignIn: just the remembered selection moved
# ---------------------------------------------------------------------------
$signIn = $wb.Worksheets.Item("SignIn")
$signIn.Range("F12").Select()

# ---------------------------------------------------------------------------
# Education: new sheet, inserted before "Profile"
# ---------------------------------------------------------------------------
$profile = $wb.Worksheets.Item("Profile")
$education = $wb.Worksheets.Add($profile)
$education.Name = "Education"

$education.Range("A1").Value = "dropDownListuni"
$education.Range("A2").Value = "JNTUH"
$education.Range("B2").Value = "Bachelor Degree"
$education.Range("B1").Value = "bachelor"

$education.Columns.Item(1).ColumnWidth = 16.67

# ---------------------------------------------------------------------------
# Language: new sheet, inserted right after Education (before Profile)
# ---------------------------------------------------------------------------
$language = $wb.Worksheets.Add($null, $education)
$language.Name = "Language"

$language.Range("A1").Value = "addLang"
$language.Range("A2").Value = "Spanish"
$language.Range("B1").Value = "languageModify"
$language.Range("B2").Value = "Telugu"

# header row styled like the other sheets' header rows (shaded fill)
$signIn.Range("A1").Copy() | Out-Null
$language.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$language.Columns.Item(1).ColumnWidth = 7.5
$language.Columns.Item(2).ColumnWidth = 14.67

# ---------------------------------------------------------------------------
# Profile: remembered selection moved
# ---------------------------------------------------------------------------
$profile.Range("D29").Select()

# ---------------------------------------------------------------------------
# Certificate: new sheet, inserted after EditShareSkill
# ---------------------------------------------------------------------------
$editShareSkill = $wb.Worksheets.Item("EditShareSkill")
$certificate = $wb.Worksheets.Add($null, $editShareSkill)
$certificate.Name = "Certificate"

$certificate.Range("A1").Value = "cert"
$certificate.Range("B1").Value = "certFrom"
$certificate.Range("C1").Value = "certificationmodify"
$certificate.Range("A2").Value = "ISTQB"
$certificate.Range("B2").Value = "ANZTB"
$certificate.Range("C2").Value = "QTP"

$certificate.Range("C2").Select()

# ---------------------------------------------------------------------------
# Four blank placeholder sheets appended at the end of the workbook
# ---------------------------------------------------------------------------
$prev = $certificate
foreach ($nm in @("Sheet4", "Sheet5", "Sheet6", "Sheet7")) {
    $blank = $wb.Worksheets.Add($null, $prev)
    $blank.Name = $nm
    $prev = $blank
}

# ---------------------------------------------------------------------------
# Education becomes the active tab
# ---------------------------------------------------------------------------
$education.Activate()
$education.Range("J8").Select()
